# Highlight three reviewer-comment paragraphs in the "PreformattedText"
# feedback list, matching the author's latest round of edits:
#   - "-sentence line 258-259: ..."                -> yellow highlight
#   - "-line 261: ..."                              -> green (bright green) highlight
#   - "-equation 7: Your equation ends with a period ..." -> green (bright green) highlight
#
# wdYellow = 7, wdBrightGreen = 4 (WdColorIndex enumeration)

$d = $word.ActiveDocument

$wdYellow = 7
$wdBrightGreen = 4

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text

    if ($text -like "-sentence lin*258-259*intersection area alone is not a probability*") {
        $para.Range.Font.HighlightColorIndex = $wdYellow
    }
    elseif ($text -like "-line 261:*integrated out (marginalized)*") {
        $para.Range.Font.HighlightColorIndex = $wdBrightGreen
    }
    elseif ($text -like "-equation 7: Your equation ends with a period*") {
        $para.Range.Font.HighlightColorIndex = $wdBrightGreen
    }
}
